$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Test Cases")

# --- New test case row 51 (WAT53) ---------------------------------------
$ws.Cells.Item(51, 1).Value = "WAT53"
# --- New test case row 52 (WAT54) ---------------------------------------
$ws.Cells.Item(52, 1).Value = "WAT54"

# Row 51 JIRA id(s) + description
$ws.Cells.Item(51, 2).Value = "WAT-372||WAT-369||WAT-636"
$ws.Cells.Item(51, 3).Value = "erify that Help file is accessible from within the application via the 'Help' link on the profile menu||Verify that system makes available the WAT help file (PDF) as a static link||Verify that user should able to access feedback survey page using 'Feedback' link on the profile menu"

# Row 52 JIRA id(s) + description
$ws.Cells.Item(52, 2).Value = "WAT-371||WAT-370"
$ws.Cells.Item(52, 3).Value = "Verify that 'Help' link available under profile flyout||Verify that entitled user for WAT have access to the help file"

# Runmode column
$ws.Cells.Item(51, 4).Value = "Y"
$ws.Cells.Item(52, 4).Value = "Y"

# Results column stays blank (same as the rest of the sheet)
$ws.Cells.Item(51, 5).Value = ""
$ws.Cells.Item(52, 5).Value = ""

# --- Formatting: match the bordered/wrapped look used by the rest of the
# table (thin black border all round, description column wraps text). ---
$rng51 = $ws.Range("A51:E51")
$rng52 = $ws.Range("A52:E52")

$rng51.Borders.ColorIndex = 1
$rng51.Borders.LineStyle = 1
$rng52.Borders.ColorIndex = 1
$rng52.Borders.LineStyle = 1

$ws.Range("C51").WrapText = $true
$ws.Range("C52").WrapText = $true

# Row 51 holds a long wrapped description, so it needs extra height.
$ws.Rows.Item(51).RowHeight = 30

# --- Selection / dimension bookkeeping ----------------------------------
$null = $ws.Range("D52").Select()
